# Apply minor numeric corrections to the coworker notebook data
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 73.3

$ws.Range("B3").Value = 28.4
$ws.Range("C3").Value = 78

$ws.Range("C4").Value = 73.90000000000001

$ws.Range("C5").Value = 40.9

$ws.Range("C6").Value = 11.7

$ws.Range("C7").Value = 20.7

$ws.Range("C10").Value = 29.3

$ws.Range("C11").Value = 88.7

$ws.Range("C14").Value = 23.9

$ws.Range("C16").Value = 29.7

$ws.Range("C17").Value = 101.3

$ws.Range("C19").Value = 117.1
